$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header D1 from "frequency_count" to "frequency"
$ws.Range("D1").Value = "frequency"

# Fill in course_id (column B) values for rows 1283-1387.
# Pattern: blocks of 21 consecutive rows share the same course_id,
# continuing the sequence already present above (... 91, 92, 93, 94, 95, 96, 97 ...)
$startRow = 1283
$courseId = 93
$row = $startRow

for ($block = 0; $block -lt 5; $block++) {
    for ($i = 0; $i -lt 21; $i++) {
        $ws.Cells.Item($row, 2).Value = $courseId
        $row = $row + 1
    }
    $courseId = $courseId + 1
}
